$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.05299063594014523
$ws.Range("E2").Value = 0.05299063594014523

$ws.Range("D3").Value = 0.000000000005106535728109117
$ws.Range("E3").Value = 0.000000000005106535728109117

$ws.Range("D4").Value = 0.000000000000000000000000000000000002342648270598597
$ws.Range("E4").Value = 0.000000000000000000000000000000000002342648270598597

$ws.Range("D5").Value = 0.000000000000000000003523455761489043
$ws.Range("E5").Value = 0.000000000000000000003523455761489043

$ws.Range("D6").Value = 0.9962727436981004
$ws.Range("E6").Value = 0.9962727436981004

$ws.Range("D7").Value = 0.0000002923365786261555
$ws.Range("E7").Value = 0.9999997076634214

$ws.Range("D9").Value = 0.892556217393733
$ws.Range("E9").Value = 0.107443782606267

$ws.Range("D10").Value = 0.0002368230540009561
$ws.Range("E10").Value = 0.999763176945999

$ws.Range("D11").Value = 0.0000000004657609661329626
$ws.Range("E11").Value = 0.999999999534239
$ws.Range("F11").Value = 5.06411075592041
